# AmongUS Assignment - add Section 3 completion (Question 3.5 rename + new Question 3.6 slide)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Edit existing slide 15 ("Question 3.4" -> "Question 3.5", new command /
#    answer text, resize answer box, add screenshot picture)
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

# -- Shape 1: title textbox: "Question 3.4 - " -> "Question 3.5 - "
$titleShape = $s15.Shapes.Item(1)
$titleRuns = $titleShape.TextFrame.TextRange.Runs()
$titleRuns.Runs(1).Text = 'Question 3.5 - '

# -- Shape 2: command placeholder - keep bold "Command –" run, replace the
#    rest of the command text with the new aggregation pipeline
$cmdShape = $s15.Shapes.Item(2)
$cmdRange = $cmdShape.TextFrame.TextRange
$cmdFull = $cmdRange.Text
$firstRunLen = 'Command –'.Length
$restLen = $cmdFull.Length - $firstRunLen
$cmdRest = $cmdRange.Characters($firstRunLen + 1, $restLen)
$cmdRest.Text = ' AMONGUS>db.Game.aggregate([{ $unwind: ''$Game_Feed'' },{$project: {''Game_Feed'':1}},{$match:{''Game_Feed.Outcome'':{ $regex: "End$"}}},{ $group: { _id: { Action: "$Game_Feed.Action" }, totalwins: { $sum:1} } }])'

# -- Shape 3: answer placeholder - resize box and update the answer text
$ansShape = $s15.Shapes.Item(3)
$ansShape.Left = 673391 / 12700
$ansShape.Top = 2748451 / 12700
$ansShape.Width = 10089107 / 12700
$ansShape.Height = 845460 / 12700

$ansRuns = $ansShape.TextFrame.TextRange.Runs()
$ansRuns.Runs(2).Text = 'The crew wins  voting 279 times'

# -- New picture: reuse the terminal-screenshot image already embedded on
#    slide 12 (same command/result, same size) by copy/pasting it in place.
$s12 = $p.Slides.Item(12)
$srcPic = $s12.Shapes.Item(4)
$srcPic.Copy()
$pastedShapes = $s15.Shapes.Paste()
$newPic = $pastedShapes.Item(1)
$newPic.Left = 445261 / 12700
$newPic.Top = 3939680 / 12700
$newPic.Width = 11496529 / 12700
$newPic.Height = 1227494 / 12700

# ---------------------------------------------------------------------------
# 2) Duplicate slide 15 -> new slide 16 ("Question 3.6")
# ---------------------------------------------------------------------------
$s15.Duplicate() | Out-Null
$s16 = $p.Slides.Item(16)

# -- Shape 1: title textbox - replace with the Question 3.6 copy
$titleShape16 = $s16.Shapes.Item(1)
$titleRuns16 = $titleShape16.TextFrame.TextRange.Runs()
$titleRuns16.Runs(1).Text = 'Question 3.6 - '
$titleRuns16.Runs(2).Text = 'Overall aggregation - '
$titleRuns16.Runs(3).Text = 'The questions you answered in this task were all related to high-level aggregations across the entire collection. In your opinion, is the game more or less hard for impostors? Justify your answer with suitable insights from the data.'
$titleShape16.Height = 923330 / 12700

# -- Shape 3: answer placeholder - replace with the final overall answer
$ansShape16 = $s16.Shapes.Item(3)
$ansRuns16 = $ansShape16.TextFrame.TextRange.Runs()
$ansRuns16.Runs(2).Text = 'Imposter won 35.27% (176 out of 499) of the time and Crew won 64.7% (323 out of 499). Since the odds of winning by Crew is ~65%, the game is hard for imposters. '

Write-Host "Done. Slide count:" $p.Slides.Count
